$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Usedcarschn Test: STARTED"
$ws.Range("A2").Value = "Usedcarschn Test Case: startBrowser Test Method: SUCCESS"
$ws.Range("A3").Value = "Usedcarschn Test Case: clickSearchUsedCarsFromDrpdwn Test Method: SUCCESS"
$ws.Range("A4").Value = "Usedcarschn Test Case: navigateToUsedCarsPage Test Method: SUCCESS"
$ws.Range("A5").Value = "Usedcarschn Test Case: outputDisplay Test Method: SUCCESS"
$ws.Range("A6").Value = "Usedcarschn Test: ENDED"

$ws.Range("A7:B8").Delete()
